# Auto-generated PowerShell Excel COM-interop script
# Applies the Wisconsin overview government-grants table restructuring:
#  - Renames 'Operating surplus ... (%)' columns to 'Size of operating surplus ...'
#  - Moves 'Share of 990 filers with government grants at risk' to be the first data column
#  - Relabels Congressional District / Size / Subsector row categories
#  - Reorders the Size and Subsector rows to their new canonical order

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Overall")
$ws.Cells.Item(1,1).Value = "Share of 990 filers with government grants at risk"
$ws.Cells.Item(1,2).Value = "Number of 990 filers with government grants"
$ws.Cells.Item(1,3).Value = "Total government grants (`$)"
$ws.Cells.Item(1,4).Value = "Size of operating surplus with government grants"
$ws.Cells.Item(1,5).Value = "Size of operating surplus without government grants"

$ws.Cells.Item(2,1).Value = "'63.53%"
$ws.Cells.Item(2,2).Value = "'2,155"
$ws.Cells.Item(2,3).Value = "'`$3,131,962,290"
$ws.Cells.Item(2,4).Value = "'10.44%"
$ws.Cells.Item(2,5).Value = "'-8.25%"


$ws = $wb.Worksheets.Item("County")
$ws.Cells.Item(1,1).Value = "Geography"
$ws.Cells.Item(1,2).Value = "Share of 990 filers with government grants at risk"
$ws.Cells.Item(1,3).Value = "Number of 990 filers with government grants"
$ws.Cells.Item(1,4).Value = "Total government grants (`$)"
$ws.Cells.Item(1,5).Value = "Size of operating surplus with government grants"
$ws.Cells.Item(1,6).Value = "Size of operating surplus without government grants"

$ws.Cells.Item(2,1).Value = "United States"
$ws.Cells.Item(2,2).Value = "'67.35%"
$ws.Cells.Item(2,3).Value = "'103,475"
$ws.Cells.Item(2,4).Value = "'`$267,700,640,005"
$ws.Cells.Item(2,5).Value = "'9.05%"
$ws.Cells.Item(2,6).Value = "'-12.83%"

$ws.Cells.Item(3,1).Value = "Wisconsin"
$ws.Cells.Item(3,2).Value = "'63.53%"
$ws.Cells.Item(3,3).Value = "'2,155"
$ws.Cells.Item(3,4).Value = "'`$3,131,962,290"
$ws.Cells.Item(3,5).Value = "'10.44%"
$ws.Cells.Item(3,6).Value = "'-8.25%"

$ws.Cells.Item(4,1).Value = "Adams County"
$ws.Cells.Item(4,2).Value = "'0.00%"
$ws.Cells.Item(4,3).Value = "'1"
$ws.Cells.Item(4,4).Value = "'`$3,926,899"
$ws.Cells.Item(4,5).Value = "'20.23%"
$ws.Cells.Item(4,6).Value = "'5.18%"

$ws.Cells.Item(5,1).Value = "Ashland County"
$ws.Cells.Item(5,2).Value = "'73.33%"
$ws.Cells.Item(5,3).Value = "'15"
$ws.Cells.Item(5,4).Value = "'`$19,971,012"
$ws.Cells.Item(5,5).Value = "'11.51%"
$ws.Cells.Item(5,6).Value = "'-22.44%"

$ws.Cells.Item(6,1).Value = "Barron County"
$ws.Cells.Item(6,2).Value = "'58.33%"
$ws.Cells.Item(6,3).Value = "'12"
$ws.Cells.Item(6,4).Value = "'`$6,511,662"
$ws.Cells.Item(6,5).Value = "'19.73%"
$ws.Cells.Item(6,6).Value = "'-17.51%"

$ws.Cells.Item(7,1).Value = "Bayfield County"
$ws.Cells.Item(7,2).Value = "'50.00%"
$ws.Cells.Item(7,3).Value = "'6"
$ws.Cells.Item(7,4).Value = "'`$15,924,578"
$ws.Cells.Item(7,5).Value = "'16.80%"
$ws.Cells.Item(7,6).Value = "'0.07%"

$ws.Cells.Item(8,1).Value = "Brown County"
$ws.Cells.Item(8,2).Value = "'58.97%"
$ws.Cells.Item(8,3).Value = "'78"
$ws.Cells.Item(8,4).Value = "'`$98,132,741"
$ws.Cells.Item(8,5).Value = "'13.40%"
$ws.Cells.Item(8,6).Value = "'-4.30%"

$ws.Cells.Item(9,1).Value = "Buffalo County"
$ws.Cells.Item(9,2).Value = "'0.00%"
$ws.Cells.Item(9,3).Value = "'1"
$ws.Cells.Item(9,4).Value = "'`$8,000"
$ws.Cells.Item(9,5).Value = "'29.25%"
$ws.Cells.Item(9,6).Value = "'18.11%"

$ws.Cells.Item(10,1).Value = "Burnett County"
$ws.Cells.Item(10,2).Value = "'100.00%"
$ws.Cells.Item(10,3).Value = "'6"
$ws.Cells.Item(10,4).Value = "'`$1,155,073"
$ws.Cells.Item(10,5).Value = "'-4.06%"
$ws.Cells.Item(10,6).Value = "'-29.90%"

$ws.Cells.Item(11,1).Value = "Calumet County"
$ws.Cells.Item(11,2).Value = "'50.00%"
$ws.Cells.Item(11,3).Value = "'4"
$ws.Cells.Item(11,4).Value = "'`$493,998"
$ws.Cells.Item(11,5).Value = "'14.63%"
$ws.Cells.Item(11,6).Value = "'5.91%"

$ws.Cells.Item(12,1).Value = "Chippewa County"
$ws.Cells.Item(12,2).Value = "'58.33%"
$ws.Cells.Item(12,3).Value = "'12"
$ws.Cells.Item(12,4).Value = "'`$3,517,264"
$ws.Cells.Item(12,5).Value = "'4.72%"
$ws.Cells.Item(12,6).Value = "'-10.29%"

$ws.Cells.Item(13,1).Value = "Clark County"
$ws.Cells.Item(13,2).Value = "'14.29%"
$ws.Cells.Item(13,3).Value = "'7"
$ws.Cells.Item(13,4).Value = "'`$366,944"
$ws.Cells.Item(13,5).Value = "'20.78%"
$ws.Cells.Item(13,6).Value = "'16.20%"

$ws.Cells.Item(14,1).Value = "Columbia County"
$ws.Cells.Item(14,2).Value = "'63.64%"
$ws.Cells.Item(14,3).Value = "'11"
$ws.Cells.Item(14,4).Value = "'`$21,978,621"
$ws.Cells.Item(14,5).Value = "'9.81%"
$ws.Cells.Item(14,6).Value = "'-5.20%"

$ws.Cells.Item(15,1).Value = "Crawford County"
$ws.Cells.Item(15,2).Value = "'50.00%"
$ws.Cells.Item(15,3).Value = "'6"
$ws.Cells.Item(15,4).Value = "'`$5,293,386"
$ws.Cells.Item(15,5).Value = "'13.92%"
$ws.Cells.Item(15,6).Value = "'1.52%"

$ws.Cells.Item(16,1).Value = "Dane County"
$ws.Cells.Item(16,2).Value = "'63.46%"
$ws.Cells.Item(16,3).Value = "'405"
$ws.Cells.Item(16,4).Value = "'`$528,520,366"
$ws.Cells.Item(16,5).Value = "'11.15%"
$ws.Cells.Item(16,6).Value = "'-7.99%"

$ws.Cells.Item(17,1).Value = "Dodge County"
$ws.Cells.Item(17,2).Value = "'66.67%"
$ws.Cells.Item(17,3).Value = "'24"
$ws.Cells.Item(17,4).Value = "'`$11,130,682"
$ws.Cells.Item(17,5).Value = "'9.01%"
$ws.Cells.Item(17,6).Value = "'-3.93%"

$ws.Cells.Item(18,1).Value = "Door County"
$ws.Cells.Item(18,2).Value = "'48.39%"
$ws.Cells.Item(18,3).Value = "'31"
$ws.Cells.Item(18,4).Value = "'`$9,027,524"
$ws.Cells.Item(18,5).Value = "'22.53%"
$ws.Cells.Item(18,6).Value = "'3.09%"

$ws.Cells.Item(19,1).Value = "Douglas County"
$ws.Cells.Item(19,2).Value = "'92.86%"
$ws.Cells.Item(19,3).Value = "'28"
$ws.Cells.Item(19,4).Value = "'`$25,241,177"
$ws.Cells.Item(19,5).Value = "'0.22%"
$ws.Cells.Item(19,6).Value = "'-55.02%"

$ws.Cells.Item(20,1).Value = "Dunn County"
$ws.Cells.Item(20,2).Value = "'75.00%"
$ws.Cells.Item(20,3).Value = "'16"
$ws.Cells.Item(20,4).Value = "'`$15,251,893"
$ws.Cells.Item(20,5).Value = "'5.17%"
$ws.Cells.Item(20,6).Value = "'-18.47%"

$ws.Cells.Item(21,1).Value = "Eau Claire County"
$ws.Cells.Item(21,2).Value = "'68.29%"
$ws.Cells.Item(21,3).Value = "'41"
$ws.Cells.Item(21,4).Value = "'`$18,468,329"
$ws.Cells.Item(21,5).Value = "'9.98%"
$ws.Cells.Item(21,6).Value = "'-6.13%"

$ws.Cells.Item(22,1).Value = "Florence County"
$ws.Cells.Item(22,2).Value = "'0.00%"
$ws.Cells.Item(22,3).Value = "'1"
$ws.Cells.Item(22,4).Value = "'`$2,020"
$ws.Cells.Item(22,5).Value = "'25.31%"
$ws.Cells.Item(22,6).Value = "'21.73%"

$ws.Cells.Item(23,1).Value = "Fond du Lac County"
$ws.Cells.Item(23,2).Value = "'58.97%"
$ws.Cells.Item(23,3).Value = "'39"
$ws.Cells.Item(23,4).Value = "'`$334,647,062"
$ws.Cells.Item(23,5).Value = "'10.59%"
$ws.Cells.Item(23,6).Value = "'-5.34%"

$ws.Cells.Item(24,1).Value = "Forest County"
$ws.Cells.Item(24,2).Value = "'0.00%"
$ws.Cells.Item(24,3).Value = "'1"
$ws.Cells.Item(24,4).Value = "'`$66,427"
$ws.Cells.Item(24,5).Value = "'10.88%"
$ws.Cells.Item(24,6).Value = "'2.22%"

$ws.Cells.Item(25,1).Value = "Grant County"
$ws.Cells.Item(25,2).Value = "'66.67%"
$ws.Cells.Item(25,3).Value = "'15"
$ws.Cells.Item(25,4).Value = "'`$15,775,279"
$ws.Cells.Item(25,5).Value = "'3.71%"
$ws.Cells.Item(25,6).Value = "'-12.01%"

$ws.Cells.Item(26,1).Value = "Green County"
$ws.Cells.Item(26,2).Value = "'80.00%"
$ws.Cells.Item(26,3).Value = "'15"
$ws.Cells.Item(26,4).Value = "'`$4,089,279"
$ws.Cells.Item(26,5).Value = "'10.46%"
$ws.Cells.Item(26,6).Value = "'-8.68%"

$ws.Cells.Item(27,1).Value = "Green Lake County"
$ws.Cells.Item(27,2).Value = "'37.50%"
$ws.Cells.Item(27,3).Value = "'8"
$ws.Cells.Item(27,4).Value = "'`$2,560,810"
$ws.Cells.Item(27,5).Value = "'31.67%"
$ws.Cells.Item(27,6).Value = "'17.49%"

$ws.Cells.Item(28,1).Value = "Iowa County"
$ws.Cells.Item(28,2).Value = "'73.68%"
$ws.Cells.Item(28,3).Value = "'19"
$ws.Cells.Item(28,4).Value = "'`$17,251,891"
$ws.Cells.Item(28,5).Value = "'16.52%"
$ws.Cells.Item(28,6).Value = "'-17.81%"

$ws.Cells.Item(29,1).Value = "Iron County"
$ws.Cells.Item(29,2).Value = "'50.00%"
$ws.Cells.Item(29,3).Value = "'2"
$ws.Cells.Item(29,4).Value = "'`$226,740"
$ws.Cells.Item(29,5).Value = "'3.37%"
$ws.Cells.Item(29,6).Value = "'-47.18%"

$ws.Cells.Item(30,1).Value = "Jackson County"
$ws.Cells.Item(30,2).Value = "'50.00%"
$ws.Cells.Item(30,3).Value = "'6"
$ws.Cells.Item(30,4).Value = "'`$6,281,889"
$ws.Cells.Item(30,5).Value = "'12.09%"
$ws.Cells.Item(30,6).Value = "'1.29%"

$ws.Cells.Item(31,1).Value = "Jefferson County"
$ws.Cells.Item(31,2).Value = "'47.06%"
$ws.Cells.Item(31,3).Value = "'17"
$ws.Cells.Item(31,4).Value = "'`$9,193,075"
$ws.Cells.Item(31,5).Value = "'20.73%"
$ws.Cells.Item(31,6).Value = "'1.29%"

$ws.Cells.Item(32,1).Value = "Juneau County"
$ws.Cells.Item(32,2).Value = "'50.00%"
$ws.Cells.Item(32,3).Value = "'4"
$ws.Cells.Item(32,4).Value = "'`$260,365"
$ws.Cells.Item(32,5).Value = "'13.89%"
$ws.Cells.Item(32,6).Value = "'-4.49%"

$ws.Cells.Item(33,1).Value = "Kenosha County"
$ws.Cells.Item(33,2).Value = "'62.07%"
$ws.Cells.Item(33,3).Value = "'29"
$ws.Cells.Item(33,4).Value = "'`$36,383,662"
$ws.Cells.Item(33,5).Value = "'11.27%"
$ws.Cells.Item(33,6).Value = "'-7.65%"

$ws.Cells.Item(34,1).Value = "Kewaunee County"
$ws.Cells.Item(34,2).Value = "'60.00%"
$ws.Cells.Item(34,3).Value = "'5"
$ws.Cells.Item(34,4).Value = "'`$186,738"
$ws.Cells.Item(34,5).Value = "'-11.82%"
$ws.Cells.Item(34,6).Value = "'-21.69%"

$ws.Cells.Item(35,1).Value = "La Crosse County"
$ws.Cells.Item(35,2).Value = "'62.96%"
$ws.Cells.Item(35,3).Value = "'54"
$ws.Cells.Item(35,4).Value = "'`$46,181,864"
$ws.Cells.Item(35,5).Value = "'7.28%"
$ws.Cells.Item(35,6).Value = "'-9.04%"

$ws.Cells.Item(36,1).Value = "Lafayette County"
$ws.Cells.Item(36,2).Value = "'66.67%"
$ws.Cells.Item(36,3).Value = "'3"
$ws.Cells.Item(36,4).Value = "'`$301,829"
$ws.Cells.Item(36,5).Value = "'9.67%"
$ws.Cells.Item(36,6).Value = "'-13.60%"

$ws.Cells.Item(37,1).Value = "Langlade County"
$ws.Cells.Item(37,2).Value = "'50.00%"
$ws.Cells.Item(37,3).Value = "'6"
$ws.Cells.Item(37,4).Value = "'`$3,411,975"
$ws.Cells.Item(37,5).Value = "'18.28%"
$ws.Cells.Item(37,6).Value = "'-4.08%"

$ws.Cells.Item(38,1).Value = "Lincoln County"
$ws.Cells.Item(38,2).Value = "'100.00%"
$ws.Cells.Item(38,3).Value = "'2"
$ws.Cells.Item(38,4).Value = "'`$507,757"
$ws.Cells.Item(38,5).Value = "'10.44%"
$ws.Cells.Item(38,6).Value = "'-28.52%"

$ws.Cells.Item(39,1).Value = "Manitowoc County"
$ws.Cells.Item(39,2).Value = "'60.71%"
$ws.Cells.Item(39,3).Value = "'28"
$ws.Cells.Item(39,4).Value = "'`$9,622,810"
$ws.Cells.Item(39,5).Value = "'10.45%"
$ws.Cells.Item(39,6).Value = "'-11.09%"

$ws.Cells.Item(40,1).Value = "Marathon County"
$ws.Cells.Item(40,2).Value = "'62.75%"
$ws.Cells.Item(40,3).Value = "'51"
$ws.Cells.Item(40,4).Value = "'`$42,555,500"
$ws.Cells.Item(40,5).Value = "'5.06%"
$ws.Cells.Item(40,6).Value = "'-8.87%"

$ws.Cells.Item(41,1).Value = "Marinette County"
$ws.Cells.Item(41,2).Value = "'57.14%"
$ws.Cells.Item(41,3).Value = "'14"
$ws.Cells.Item(41,4).Value = "'`$2,014,625"
$ws.Cells.Item(41,5).Value = "'12.20%"
$ws.Cells.Item(41,6).Value = "'-7.61%"

$ws.Cells.Item(42,1).Value = "Menominee County"
$ws.Cells.Item(42,2).Value = "'100.00%"
$ws.Cells.Item(42,3).Value = "'2"
$ws.Cells.Item(42,4).Value = "'`$10,968,478"
$ws.Cells.Item(42,5).Value = "'18.27%"
$ws.Cells.Item(42,6).Value = "'-70.22%"

$ws.Cells.Item(43,1).Value = "Milwaukee County"
$ws.Cells.Item(43,2).Value = "'70.75%"
$ws.Cells.Item(43,3).Value = "'424"
$ws.Cells.Item(43,4).Value = "'`$1,068,151,542"
$ws.Cells.Item(43,5).Value = "'7.74%"
$ws.Cells.Item(43,6).Value = "'-16.17%"

$ws.Cells.Item(44,1).Value = "Monroe County"
$ws.Cells.Item(44,2).Value = "'69.23%"
$ws.Cells.Item(44,3).Value = "'13"
$ws.Cells.Item(44,4).Value = "'`$18,564,884"
$ws.Cells.Item(44,5).Value = "'8.87%"
$ws.Cells.Item(44,6).Value = "'-3.43%"

$ws.Cells.Item(45,1).Value = "Oconto County"
$ws.Cells.Item(45,2).Value = "'75.00%"
$ws.Cells.Item(45,3).Value = "'8"
$ws.Cells.Item(45,4).Value = "'`$32,560,723"
$ws.Cells.Item(45,5).Value = "'5.08%"
$ws.Cells.Item(45,6).Value = "'-19.99%"

$ws.Cells.Item(46,1).Value = "Oneida County"
$ws.Cells.Item(46,2).Value = "'70.59%"
$ws.Cells.Item(46,3).Value = "'17"
$ws.Cells.Item(46,4).Value = "'`$7,652,190"
$ws.Cells.Item(46,5).Value = "'7.54%"
$ws.Cells.Item(46,6).Value = "'-8.70%"

$ws.Cells.Item(47,1).Value = "Outagamie County"
$ws.Cells.Item(47,2).Value = "'56.25%"
$ws.Cells.Item(47,3).Value = "'64"
$ws.Cells.Item(47,4).Value = "'`$57,685,330"
$ws.Cells.Item(47,5).Value = "'9.59%"
$ws.Cells.Item(47,6).Value = "'-3.80%"

$ws.Cells.Item(48,1).Value = "Ozaukee County"
$ws.Cells.Item(48,2).Value = "'64.52%"
$ws.Cells.Item(48,3).Value = "'31"
$ws.Cells.Item(48,4).Value = "'`$25,808,366"
$ws.Cells.Item(48,5).Value = "'12.18%"
$ws.Cells.Item(48,6).Value = "'-7.42%"

$ws.Cells.Item(49,1).Value = "Pepin County"
$ws.Cells.Item(49,2).Value = "'0.00%"
$ws.Cells.Item(49,3).Value = "'1"
$ws.Cells.Item(49,4).Value = "'`$218,178"
$ws.Cells.Item(49,5).Value = "'5.44%"
$ws.Cells.Item(49,6).Value = "'4.22%"

$ws.Cells.Item(50,1).Value = "Pierce County"
$ws.Cells.Item(50,2).Value = "'71.43%"
$ws.Cells.Item(50,3).Value = "'7"
$ws.Cells.Item(50,4).Value = "'`$1,960,910"
$ws.Cells.Item(50,5).Value = "'10.87%"
$ws.Cells.Item(50,6).Value = "'-7.86%"

$ws.Cells.Item(51,1).Value = "Polk County"
$ws.Cells.Item(51,2).Value = "'66.67%"
$ws.Cells.Item(51,3).Value = "'12"
$ws.Cells.Item(51,4).Value = "'`$10,106,572"
$ws.Cells.Item(51,5).Value = "'5.29%"
$ws.Cells.Item(51,6).Value = "'-7.58%"

$ws.Cells.Item(52,1).Value = "Portage County"
$ws.Cells.Item(52,2).Value = "'57.69%"
$ws.Cells.Item(52,3).Value = "'26"
$ws.Cells.Item(52,4).Value = "'`$21,718,654"
$ws.Cells.Item(52,5).Value = "'11.65%"
$ws.Cells.Item(52,6).Value = "'-2.23%"

$ws.Cells.Item(53,1).Value = "Price County"
$ws.Cells.Item(53,2).Value = "'0.00%"
$ws.Cells.Item(53,3).Value = "'1"
$ws.Cells.Item(53,4).Value = "'`$546,408"
$ws.Cells.Item(53,5).Value = "'41.99%"
$ws.Cells.Item(53,6).Value = "'40.49%"

$ws.Cells.Item(54,1).Value = "Racine County"
$ws.Cells.Item(54,2).Value = "'70.21%"
$ws.Cells.Item(54,3).Value = "'47"
$ws.Cells.Item(54,4).Value = "'`$35,276,634"
$ws.Cells.Item(54,5).Value = "'8.52%"
$ws.Cells.Item(54,6).Value = "'-17.24%"

$ws.Cells.Item(55,1).Value = "Richland County"
$ws.Cells.Item(55,2).Value = "'50.00%"
$ws.Cells.Item(55,3).Value = "'4"
$ws.Cells.Item(55,4).Value = "'`$6,808,413"
$ws.Cells.Item(55,5).Value = "'10.82%"
$ws.Cells.Item(55,6).Value = "'-26.11%"

$ws.Cells.Item(56,1).Value = "Rock County"
$ws.Cells.Item(56,2).Value = "'56.00%"
$ws.Cells.Item(56,3).Value = "'50"
$ws.Cells.Item(56,4).Value = "'`$52,947,247"
$ws.Cells.Item(56,5).Value = "'15.55%"
$ws.Cells.Item(56,6).Value = "'-5.42%"

$ws.Cells.Item(57,1).Value = "Rusk County"
$ws.Cells.Item(57,2).Value = "'100.00%"
$ws.Cells.Item(57,3).Value = "'2"
$ws.Cells.Item(57,4).Value = "'`$6,069,065"
$ws.Cells.Item(57,5).Value = "'-4.31%"
$ws.Cells.Item(57,6).Value = "'-83.61%"

$ws.Cells.Item(58,1).Value = "Sauk County"
$ws.Cells.Item(58,2).Value = "'41.18%"
$ws.Cells.Item(58,3).Value = "'17"
$ws.Cells.Item(58,4).Value = "'`$6,016,362"
$ws.Cells.Item(58,5).Value = "'7.59%"
$ws.Cells.Item(58,6).Value = "'3.03%"

$ws.Cells.Item(59,1).Value = "Sawyer County"
$ws.Cells.Item(59,2).Value = "'58.82%"
$ws.Cells.Item(59,3).Value = "'17"
$ws.Cells.Item(59,4).Value = "'`$29,615,889"
$ws.Cells.Item(59,5).Value = "'13.34%"
$ws.Cells.Item(59,6).Value = "'-31.82%"

$ws.Cells.Item(60,1).Value = "Shawano County"
$ws.Cells.Item(60,2).Value = "'87.50%"
$ws.Cells.Item(60,3).Value = "'8"
$ws.Cells.Item(60,4).Value = "'`$1,607,319"
$ws.Cells.Item(60,5).Value = "'-1.47%"
$ws.Cells.Item(60,6).Value = "'-29.43%"

$ws.Cells.Item(61,1).Value = "Sheboygan County"
$ws.Cells.Item(61,2).Value = "'55.10%"
$ws.Cells.Item(61,3).Value = "'49"
$ws.Cells.Item(61,4).Value = "'`$29,927,511"
$ws.Cells.Item(61,5).Value = "'14.06%"
$ws.Cells.Item(61,6).Value = "'-2.44%"

$ws.Cells.Item(62,1).Value = "St. Croix County"
$ws.Cells.Item(62,2).Value = "'59.09%"
$ws.Cells.Item(62,3).Value = "'22"
$ws.Cells.Item(62,4).Value = "'`$19,100,863"
$ws.Cells.Item(62,5).Value = "'11.54%"
$ws.Cells.Item(62,6).Value = "'-1.86%"

$ws.Cells.Item(63,1).Value = "Taylor County"
$ws.Cells.Item(63,2).Value = "'44.44%"
$ws.Cells.Item(63,3).Value = "'9"
$ws.Cells.Item(63,4).Value = "'`$4,943,555"
$ws.Cells.Item(63,5).Value = "'16.89%"
$ws.Cells.Item(63,6).Value = "'5.65%"

$ws.Cells.Item(64,1).Value = "Trempealeau County"
$ws.Cells.Item(64,2).Value = "'66.67%"
$ws.Cells.Item(64,3).Value = "'3"
$ws.Cells.Item(64,4).Value = "'`$19,557,420"
$ws.Cells.Item(64,5).Value = "'2.93%"
$ws.Cells.Item(64,6).Value = "'-2.65%"

$ws.Cells.Item(65,1).Value = "Vernon County"
$ws.Cells.Item(65,2).Value = "'70.59%"
$ws.Cells.Item(65,3).Value = "'17"
$ws.Cells.Item(65,4).Value = "'`$10,796,320"
$ws.Cells.Item(65,5).Value = "'7.45%"
$ws.Cells.Item(65,6).Value = "'-6.35%"

$ws.Cells.Item(66,1).Value = "Vilas County"
$ws.Cells.Item(66,2).Value = "'33.33%"
$ws.Cells.Item(66,3).Value = "'6"
$ws.Cells.Item(66,4).Value = "'`$945,050"
$ws.Cells.Item(66,5).Value = "'32.38%"
$ws.Cells.Item(66,6).Value = "'15.77%"

$ws.Cells.Item(67,1).Value = "Walworth County"
$ws.Cells.Item(67,2).Value = "'58.33%"
$ws.Cells.Item(67,3).Value = "'24"
$ws.Cells.Item(67,4).Value = "'`$5,927,542"
$ws.Cells.Item(67,5).Value = "'5.14%"
$ws.Cells.Item(67,6).Value = "'-12.37%"

$ws.Cells.Item(68,1).Value = "Washburn County"
$ws.Cells.Item(68,2).Value = "'66.67%"
$ws.Cells.Item(68,3).Value = "'9"
$ws.Cells.Item(68,4).Value = "'`$2,239,318"
$ws.Cells.Item(68,5).Value = "'0.16%"
$ws.Cells.Item(68,6).Value = "'-44.65%"

$ws.Cells.Item(69,1).Value = "Washington County"
$ws.Cells.Item(69,2).Value = "'60.00%"
$ws.Cells.Item(69,3).Value = "'20"
$ws.Cells.Item(69,4).Value = "'`$11,020,704"
$ws.Cells.Item(69,5).Value = "'12.47%"
$ws.Cells.Item(69,6).Value = "'-5.75%"

$ws.Cells.Item(70,1).Value = "Waukesha County"
$ws.Cells.Item(70,2).Value = "'63.79%"
$ws.Cells.Item(70,3).Value = "'116"
$ws.Cells.Item(70,4).Value = "'`$95,298,203"
$ws.Cells.Item(70,5).Value = "'9.81%"
$ws.Cells.Item(70,6).Value = "'-8.40%"

$ws.Cells.Item(71,1).Value = "Waupaca County"
$ws.Cells.Item(71,2).Value = "'57.89%"
$ws.Cells.Item(71,3).Value = "'19"
$ws.Cells.Item(71,4).Value = "'`$2,887,397"
$ws.Cells.Item(71,5).Value = "'13.71%"
$ws.Cells.Item(71,6).Value = "'-3.55%"

$ws.Cells.Item(72,1).Value = "Waushara County"
$ws.Cells.Item(72,2).Value = "'71.43%"
$ws.Cells.Item(72,3).Value = "'7"
$ws.Cells.Item(72,4).Value = "'`$8,673,092"
$ws.Cells.Item(72,5).Value = "'-1.69%"
$ws.Cells.Item(72,6).Value = "'-14.21%"

$ws.Cells.Item(73,1).Value = "Winnebago County"
$ws.Cells.Item(73,2).Value = "'48.15%"
$ws.Cells.Item(73,3).Value = "'54"
$ws.Cells.Item(73,4).Value = "'`$84,991,840"
$ws.Cells.Item(73,5).Value = "'14.09%"
$ws.Cells.Item(73,6).Value = "'1.06%"

$ws.Cells.Item(74,1).Value = "Wood County"
$ws.Cells.Item(74,2).Value = "'65.38%"
$ws.Cells.Item(74,3).Value = "'26"
$ws.Cells.Item(74,4).Value = "'`$124,928,565"
$ws.Cells.Item(74,5).Value = "'6.51%"
$ws.Cells.Item(74,6).Value = "'-7.30%"


$ws = $wb.Worksheets.Item("Congressional District")
$ws.Cells.Item(1,1).Value = "Geography"
$ws.Cells.Item(1,2).Value = "Share of 990 filers with government grants at risk"
$ws.Cells.Item(1,3).Value = "Number of 990 filers with government grants"
$ws.Cells.Item(1,4).Value = "Total government grants (`$)"
$ws.Cells.Item(1,5).Value = "Size of operating surplus with government grants"
$ws.Cells.Item(1,6).Value = "Size of operating surplus without government grants"

$ws.Cells.Item(2,1).Value = "United States"
$ws.Cells.Item(2,2).Value = "'67.35%"
$ws.Cells.Item(2,3).Value = "'103,475"
$ws.Cells.Item(2,4).Value = "'`$267,700,640,005"
$ws.Cells.Item(2,5).Value = "'9.05%"
$ws.Cells.Item(2,6).Value = "'-12.83%"

$ws.Cells.Item(3,1).Value = "Wisconsin"
$ws.Cells.Item(3,2).Value = "'63.53%"
$ws.Cells.Item(3,3).Value = "'2,155"
$ws.Cells.Item(3,4).Value = "'`$3,131,962,290"
$ws.Cells.Item(3,5).Value = "'10.44%"
$ws.Cells.Item(3,6).Value = "'-8.25%"

$ws.Cells.Item(4,1).Value = "Congressional District 1"
$ws.Cells.Item(4,2).Value = "'61.07%"
$ws.Cells.Item(4,3).Value = "'149"
$ws.Cells.Item(4,4).Value = "'`$127,768,587"
$ws.Cells.Item(4,5).Value = "'12.71%"
$ws.Cells.Item(4,6).Value = "'-9.93%"

$ws.Cells.Item(5,1).Value = "Congressional District 2"
$ws.Cells.Item(5,2).Value = "'63.25%"
$ws.Cells.Item(5,3).Value = "'468"
$ws.Cells.Item(5,4).Value = "'`$564,217,822"
$ws.Cells.Item(5,5).Value = "'11.06%"
$ws.Cells.Item(5,6).Value = "'-7.92%"

$ws.Cells.Item(6,1).Value = "Congressional District 3"
$ws.Cells.Item(6,2).Value = "'63.18%"
$ws.Cells.Item(6,3).Value = "'239"
$ws.Cells.Item(6,4).Value = "'`$213,856,474"
$ws.Cells.Item(6,5).Value = "'9.11%"
$ws.Cells.Item(6,6).Value = "'-6.58%"

$ws.Cells.Item(7,1).Value = "Congressional District 4"
$ws.Cells.Item(7,2).Value = "'72.11%"
$ws.Cells.Item(7,3).Value = "'398"
$ws.Cells.Item(7,4).Value = "'`$1,057,965,007"
$ws.Cells.Item(7,5).Value = "'7.53%"
$ws.Cells.Item(7,6).Value = "'-17.63%"

$ws.Cells.Item(8,1).Value = "Congressional District 5"
$ws.Cells.Item(8,2).Value = "'60.56%"
$ws.Cells.Item(8,3).Value = "'180"
$ws.Cells.Item(8,4).Value = "'`$124,621,537"
$ws.Cells.Item(8,5).Value = "'10.89%"
$ws.Cells.Item(8,6).Value = "'-5.90%"

$ws.Cells.Item(9,1).Value = "Congressional District 6"
$ws.Cells.Item(9,2).Value = "'57.61%"
$ws.Cells.Item(9,3).Value = "'243"
$ws.Cells.Item(9,4).Value = "'`$525,169,172"
$ws.Cells.Item(9,5).Value = "'11.84%"
$ws.Cells.Item(9,6).Value = "'-4.76%"

$ws.Cells.Item(10,1).Value = "Congressional District 7"
$ws.Cells.Item(10,2).Value = "'64.63%"
$ws.Cells.Item(10,3).Value = "'246"
$ws.Cells.Item(10,4).Value = "'`$302,821,813"
$ws.Cells.Item(10,5).Value = "'8.65%"
$ws.Cells.Item(10,6).Value = "'-11.09%"

$ws.Cells.Item(11,1).Value = "Congressional District 8"
$ws.Cells.Item(11,2).Value = "'58.62%"
$ws.Cells.Item(11,3).Value = "'232"
$ws.Cells.Item(11,4).Value = "'`$215,541,878"
$ws.Cells.Item(11,5).Value = "'12.66%"
$ws.Cells.Item(11,6).Value = "'-4.03%"


$ws = $wb.Worksheets.Item("Size")
$ws.Cells.Item(1,1).Value = "Size"
$ws.Cells.Item(1,2).Value = "Share of 990 filers with government grants at risk"
$ws.Cells.Item(1,3).Value = "Number of 990 filers with government grants"
$ws.Cells.Item(1,4).Value = "Total government grants (`$)"
$ws.Cells.Item(1,5).Value = "Size of operating surplus with government grants"
$ws.Cells.Item(1,6).Value = "Size of operating surplus without government grants"

$ws.Cells.Item(2,1).Value = "Between `$100K and `$499K"
$ws.Cells.Item(2,2).Value = "'62.14%"
$ws.Cells.Item(2,3).Value = "'700"
$ws.Cells.Item(2,4).Value = "'`$62,299,839"
$ws.Cells.Item(2,5).Value = "'12.11%"
$ws.Cells.Item(2,6).Value = "'-8.03%"

$ws.Cells.Item(3,1).Value = "Between `$1M and `$4.99M"
$ws.Cells.Item(3,2).Value = "'65.01%"
$ws.Cells.Item(3,3).Value = "'603"
$ws.Cells.Item(3,4).Value = "'`$452,014,272"
$ws.Cells.Item(3,5).Value = "'9.58%"
$ws.Cells.Item(3,6).Value = "'-9.78%"

$ws.Cells.Item(4,1).Value = "Between `$500K and `$999K"
$ws.Cells.Item(4,2).Value = "'64.27%"
$ws.Cells.Item(4,3).Value = "'389"
$ws.Cells.Item(4,4).Value = "'`$91,162,861"
$ws.Cells.Item(4,5).Value = "'12.18%"
$ws.Cells.Item(4,6).Value = "'-9.99%"

$ws.Cells.Item(5,1).Value = "Between `$5M and `$9.99M"
$ws.Cells.Item(5,2).Value = "'67.61%"
$ws.Cells.Item(5,3).Value = "'142"
$ws.Cells.Item(5,4).Value = "'`$338,654,196"
$ws.Cells.Item(5,5).Value = "'7.11%"
$ws.Cells.Item(5,6).Value = "'-8.87%"

$ws.Cells.Item(6,1).Value = "Greater than `$10M"
$ws.Cells.Item(6,2).Value = "'61.25%"
$ws.Cells.Item(6,3).Value = "'240"
$ws.Cells.Item(6,4).Value = "'`$2,182,411,739"
$ws.Cells.Item(6,5).Value = "'6.57%"
$ws.Cells.Item(6,6).Value = "'-5.44%"

$ws.Cells.Item(7,1).Value = "Less than `$100K"
$ws.Cells.Item(7,2).Value = "'60.49%"
$ws.Cells.Item(7,3).Value = "'81"
$ws.Cells.Item(7,4).Value = "'`$5,419,383"
$ws.Cells.Item(7,5).Value = "'18.44%"
$ws.Cells.Item(7,6).Value = "'-15.76%"

$ws.Cells.Item(8,1).Value = "Total"
$ws.Cells.Item(8,2).Value = "'63.53%"
$ws.Cells.Item(8,3).Value = "'2,155"
$ws.Cells.Item(8,4).Value = "'`$3,131,962,290"
$ws.Cells.Item(8,5).Value = "'10.44%"
$ws.Cells.Item(8,6).Value = "'-8.25%"


$ws = $wb.Worksheets.Item("Subsector")
$ws.Cells.Item(1,1).Value = "Subsector"
$ws.Cells.Item(1,2).Value = "Share of 990 filers with government grants at risk"
$ws.Cells.Item(1,3).Value = "Number of 990 filers with government grants"
$ws.Cells.Item(1,4).Value = "Total government grants (`$)"
$ws.Cells.Item(1,5).Value = "Size of operating surplus with government grants"
$ws.Cells.Item(1,6).Value = "Size of operating surplus without government grants"

$ws.Cells.Item(2,1).Value = "Arts, Culture, and Humanities"
$ws.Cells.Item(2,2).Value = "'60.56%"
$ws.Cells.Item(2,3).Value = "'180"
$ws.Cells.Item(2,4).Value = "'`$87,932,277"
$ws.Cells.Item(2,5).Value = "'17.05%"
$ws.Cells.Item(2,6).Value = "'-6.12%"

$ws.Cells.Item(3,1).Value = "Education (Excluding Universities)"
$ws.Cells.Item(3,2).Value = "'63.74%"
$ws.Cells.Item(3,3).Value = "'182"
$ws.Cells.Item(3,4).Value = "'`$332,038,198"
$ws.Cells.Item(3,5).Value = "'10.29%"
$ws.Cells.Item(3,6).Value = "'-8.31%"

$ws.Cells.Item(4,1).Value = "Environment and Animals"
$ws.Cells.Item(4,2).Value = "'45.16%"
$ws.Cells.Item(4,3).Value = "'124"
$ws.Cells.Item(4,4).Value = "'`$52,126,882"
$ws.Cells.Item(4,5).Value = "'19.66%"
$ws.Cells.Item(4,6).Value = "'3.08%"

$ws.Cells.Item(5,1).Value = "Health (Excluding Hospitals)"
$ws.Cells.Item(5,2).Value = "'65.20%"
$ws.Cells.Item(5,3).Value = "'204"
$ws.Cells.Item(5,4).Value = "'`$527,687,775"
$ws.Cells.Item(5,5).Value = "'7.58%"
$ws.Cells.Item(5,6).Value = "'-10.33%"

$ws.Cells.Item(6,1).Value = "Hospitals"
$ws.Cells.Item(6,2).Value = "'40.00%"
$ws.Cells.Item(6,3).Value = "'30"
$ws.Cells.Item(6,4).Value = "'`$88,661,189"
$ws.Cells.Item(6,5).Value = "'9.73%"
$ws.Cells.Item(6,6).Value = "'2.20%"

$ws.Cells.Item(7,1).Value = "Human Services"
$ws.Cells.Item(7,2).Value = "'65.72%"
$ws.Cells.Item(7,3).Value = "'671"
$ws.Cells.Item(7,4).Value = "'`$739,275,028"
$ws.Cells.Item(7,5).Value = "'10.00%"
$ws.Cells.Item(7,6).Value = "'-11.00%"

$ws.Cells.Item(8,1).Value = "International, Foreign Affairs"
$ws.Cells.Item(8,2).Value = "'54.55%"
$ws.Cells.Item(8,3).Value = "'11"
$ws.Cells.Item(8,4).Value = "'`$5,048,554"
$ws.Cells.Item(8,5).Value = "'13.16%"
$ws.Cells.Item(8,6).Value = "'-5.80%"

$ws.Cells.Item(9,1).Value = "Mutual/Membership Benefit"
$ws.Cells.Item(9,2).Value = "'100.00%"
$ws.Cells.Item(9,3).Value = "'2"
$ws.Cells.Item(9,4).Value = "'`$8,602,845"
$ws.Cells.Item(9,5).Value = "'-7.17%"
$ws.Cells.Item(9,6).Value = "'-100.42%"

$ws.Cells.Item(10,1).Value = "Public, Societal Benefit"
$ws.Cells.Item(10,2).Value = "'64.13%"
$ws.Cells.Item(10,3).Value = "'184"
$ws.Cells.Item(10,4).Value = "'`$215,512,430"
$ws.Cells.Item(10,5).Value = "'12.72%"
$ws.Cells.Item(10,6).Value = "'-12.07%"

$ws.Cells.Item(11,1).Value = "Religion Related"
$ws.Cells.Item(11,2).Value = "'47.83%"
$ws.Cells.Item(11,3).Value = "'46"
$ws.Cells.Item(11,4).Value = "'`$11,287,315"
$ws.Cells.Item(11,5).Value = "'10.29%"
$ws.Cells.Item(11,6).Value = "'0.61%"

$ws.Cells.Item(12,1).Value = "Unclassified"
$ws.Cells.Item(12,2).Value = "'68.55%"
$ws.Cells.Item(12,3).Value = "'496"
$ws.Cells.Item(12,4).Value = "'`$733,082,122"
$ws.Cells.Item(12,5).Value = "'8.38%"
$ws.Cells.Item(12,6).Value = "'-10.32%"

$ws.Cells.Item(13,1).Value = "Universities"
$ws.Cells.Item(13,2).Value = "'56.00%"
$ws.Cells.Item(13,3).Value = "'25"
$ws.Cells.Item(13,4).Value = "'`$330,707,675"
$ws.Cells.Item(13,5).Value = "'6.28%"
$ws.Cells.Item(13,6).Value = "'-1.08%"

$ws.Cells.Item(14,1).Value = "Total"
$ws.Cells.Item(14,2).Value = "'63.53%"
$ws.Cells.Item(14,3).Value = "'2,155"
$ws.Cells.Item(14,4).Value = "'`$3,131,962,290"
$ws.Cells.Item(14,5).Value = "'10.44%"
$ws.Cells.Item(14,6).Value = "'-8.25%"

